# Drop the 9 oldest accelerometer samples (rows 2-10), shifting the remaining rows up,
# then append the new samples collected on may 9th.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data up by deleting the oldest 9 rows of samples
$ws.Range("A2:C10").Delete(-4162)

# Append the newly collected samples
$newSamples = @(
    @(-8.376615575381617, -8.225071150490136, -15.93280724116742),
    @(0.8738524981906739, -4.53215429782872, 15.28993327277035),
    @(-0.01965822492311098, 16.21388537543142, 17.33802195957751),
    @(-10.3103993143354, -31.18352238791278, 18.67163455826883),
    @(22.02176739828893, -22.72049116407177, -2.185512433732939),
    @(5.463320446014379, 6.157036137580864, -27.94153518676758),
    @(-37.06017020089279, -23.85922709873728, 0.4766027178080634),
    @(-19.13142990640274, -10.62035506112288, 8.715779474803668),
    @(6.810042284216134, 8.792214941978534, -5.953527409689759),
    @(-0.03172696658525354, -6.182636371680525, -11.13081671851015),
    @(12.65357358115054, -1.313164901733439, -15.5251411982945),
    @(8.714159175328177, -9.301818609237403, -17.37430092947822),
    @(-5.424608945846559, -50.41501617431641, -2.802897453308105),
    @(3.029158597333116, -28.14689919608023, 5.423780044487627),
    @(23.23876049859135, 5.670972211020553, 3.751934364863664),
    @(12.21474557604104, -7.668160905156817, -2.514941801343662),
    @(1.01306713649177, 2.441266550336438, -15.73728098188129),
    @(-0.6987579890658431, -8.856165681566392, -19.26768711635043),
    @(-16.89904033115963, -34.02156884329665, 7.39472787039573)
)

for ($i = 0; $i -lt $newSamples.Length; $i++) {
    $row = 13 + $i
    $ws.Cells.Item($row, 1).Value = $newSamples[$i][0]
    $ws.Cells.Item($row, 2).Value = $newSamples[$i][1]
    $ws.Cells.Item($row, 3).Value = $newSamples[$i][2]
}
